# "modificacion de reporte serenity"
# Hoja1: remove the two sample hyperlink rows' email-specific data (email,
# placa/address extras, phone) down to a mostly-blank template, keep the
# dropdown-driven "NA" placeholders, and drop the mailto: hyperlinks.
# Hoja2 and Hoja1 selections/scroll position are updated as well.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Hoja1 --------------------------------------------------------------

# Drop both mailto: hyperlinks (on R2 / R3) entirely.
$ws1.Hyperlinks.Delete()

# Row 2: keep the "NA" dropdown placeholders (B2:L2, P2:Q2), the address /
# phone-type / phone-number sample (M2:O2), but clear the Placa (A2) and
# the former hyperlink + cellphone cells (R2:S2).
$ws1.Range("A2").ClearContents()
$ws1.Range("R2:S2").ClearContents()

# Row 3 becomes an (almost) empty template row: only A3 / N3 / R3 keep
# their formatting (Placa style, phone-type style, hyperlink style) but
# all values are cleared.
$ws1.Range("A3:S3").ClearContents()

# Sheet1 selection / scroll position.
$ws1.Activate() | Out-Null
$ws1.Range("R3").Select() | Out-Null
$excel.ActiveWindow.ScrollColumn = 7

# New column width for column R (18), matching column L's best-fit width.
$ws1.Columns.Item(18).ColumnWidth = $ws1.Columns.Item(12).ColumnWidth

# --- Hoja2 ----------------------------------------------------------------

$ws2.Range("B9").Select() | Out-Null

# Re-activate Hoja1 so it remains the selected/active tab, matching the
# original workbook (only Hoja1 had tabSelected).
$ws1.Activate() | Out-Null
$ws1.Range("R3").Select() | Out-Null
